# Insert a new weekly record row at row 154 for
# "Terminal Hortofrutícola Agro Chillán" / Cebollín, shifting the
# existing rows 154:181 down to 155:182.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(154).Insert()

$ws.Cells.Item(154, 1).Value = 7
$ws.Cells.Item(154, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(154, 3).Value = "Ñuble"
$ws.Cells.Item(154, 4).Value = "7/20/2023"
$ws.Cells.Item(154, 5).Value = 16
$ws.Cells.Item(154, 6).Value = 100112037
$ws.Cells.Item(154, 7).Value = "Cebollín"
$ws.Cells.Item(154, 8).Value = "Sin especificar"
$ws.Cells.Item(154, 9).Value = "Segunda"
$ws.Cells.Item(154, 10).Value = 100
$ws.Cells.Item(154, 11).Value = 6000
$ws.Cells.Item(154, 12).Value = 6000
$ws.Cells.Item(154, 13).Value = 6000
$ws.Cells.Item(154, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(154, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(154, 16).Value = 167
$ws.Cells.Item(154, 17).Value = 36
$ws.Cells.Item(154, 18).Value = "Hortaliza"
